# Regenerate s_vals data to filter save games.
# Updates the numeric stat columns (B:G) for rows 2-10 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(0.1554434735375247, 1.65323645889881,   0.1529057820181812, 6.48142807727062,  0, 8.443013791725136)
    3  = @(1.505614041169197,  1.65323645889881,   0.1529057820181812, 6.48142807727062,  0, 9.793184359356808)
    4  = @(0.1554434735375247, 1.65323645889881,   0.1529057820181812, 0.4998867070740569, 0, 2.461472421528573)
    5  = @(3.182878228561681,  1.65323645889881,   0.7127328510149897, 246.9852506941017, 1, 252.5340982325772)
    6  = @(1.505614041169197,  0.3375848360084654, 0.1529057820181812, 0.4998867070740569, 0, 2.495991366269901)
    7  = @(0.3464964993005633, 1766.335244827366,  157.8057217802531,  71517.89157740913, 0, 73442.37904051605)
    8  = @(0.02258322285507441, 0.3375848360084654, 0.7127328510149897, 6.48142807727062,  0, 7.55432898714915)
    9  = @(0.1554434735375247, 1.65323645889881,   157.8057217802531,  6.48142807727062,  0, 166.0958297899601)
    10 = @(0.1554434735375247, 0.05231270169004087, 0.7127328510149897, 0.4998867070740569, 1, 1.420375733316612)
}

$columns = @("B", "C", "D", "E", "F", "G")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $columns.Length; $i++) {
        $ws.Range("$($columns[$i])$row").Value = $values[$i]
    }
}
